$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set value "O" in E6, E7, E8 (same value as already present in B6:D6, B7:D7, B8:D8)
$ws.Range("E6").Value = "O"
$ws.Range("E7").Value = "O"
$ws.Range("E8").Value = "O"

# Update selection to E9 to match the final cursor position
$ws.Range("E9").Select()
